$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 4 ----
$ws.Range("A4").Value = "Jino"
$ws.Range("B4").Value = "Mukesh"
$ws.Range("C4").Value = "Kazhakuttom"
$ws.Range("G4").Value = "Admin"
$ws.Range("E4").Value = 9876543210
$ws.Range("F2").Copy()
$ws.Range("F4").PasteSpecial(-4122)
$ws.Range("F4").Value = 32724

# ---- Row 5 ----
$ws.Range("A5").Value = "John"
$ws.Range("B5").Value = "Dhas"
$ws.Range("C5").Value = "Chennai"
$ws.Range("D5").Value = "john@gmail.com"
$ws.Range("G5").Value = "Driver"
$ws.Range("E5").Value = 784596520
$ws.Range("F2").Copy()
$ws.Range("F5").PasteSpecial(-4122)
$ws.Range("F5").Value = 32697

# ---- Row 6 ----
$ws.Range("A6").Value = "Jaiden"
$ws.Range("B6").Value = "MB"
$ws.Range("C6").Value = "Marthandam"
$ws.Range("D6").Value = "jai@yourdomain.com"
$ws.Range("G6").Value = "Cleaner"
$ws.Range("E6").Value = 7845129630
$ws.Range("F2").Copy()
$ws.Range("F6").PasteSpecial(-4122)
$ws.Range("F6").Value = 42590

# D4 written last among the new string cells so shared-string order matches
$ws.Range("D4").Value = "mukesh@gmail.com"

$excel.CutCopyMode = $false

# ---- Hyperlinks for the new email cells ----
$ws.Hyperlinks.Add($ws.Range("D4"), "mailto:mukesh@gmail.com")
$ws.Range("D2").Copy()
$ws.Range("D4").PasteSpecial(-4122)

$ws.Hyperlinks.Add($ws.Range("D5"), "mailto:john@gmail.com")
$ws.Range("D2").Copy()
$ws.Range("D5").PasteSpecial(-4122)

$ws.Hyperlinks.Add($ws.Range("D6"), "mailto:jai@yourdomain.com")
$ws.Range("D2").Copy()
$ws.Range("D6").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# ---- Selection change ----
$ws.Range("D5").Select()
